$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task names for rows 24-28 (Withdraw feature tasks) and add new
# "Real time" values for the first two of them.
$ws.Range("B24").Value = "Dodavanje Deposit metode na servis za komunikaciju sa bankom"
$ws.Range("D24").Value = 5

$ws.Range("B25").Value = "Dodavanje Withdraw metode na WalletService"
$ws.Range("D25").Value = 5

$ws.Range("B26").Value = "Implementacija testova za Withdraw"

$ws.Range("B27").Value = "Dodavanje rute za Withdraw sredstava u WalletController"

$ws.Range("B28").Value = "Dodavanje stranice za Withdraw sredstava u MVC aplikaciju"

# Move the selection like the author did while finishing up edits.
$ws.Range("B29").Select()
